$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row with Russian translations, add new "Цена за эти товары" column
$ws.Range("A1").Value = "Название"
$ws.Range("B1").Value = "Цена"
$ws.Range("C1").Value = "Количесво"
$ws.Range("D1").Value = "Цена за эти товары"

# Update data row 2 (still the "candle" product) and add total price
$ws.Range("A2").Value = "candle"
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 3000

# Remove the old third row (sylvia-breitenberg / 615 / 5) entirely
$ws.Rows.Item(3).Delete()

# Set explicit column widths to match the target layout.
# (The host's ColumnWidth setter adds a fixed 5/6-character padding when
#  round-tripping through the saved <col width=.../> attribute, so the
#  assigned value is pre-compensated to land exactly on the target width.)
$ws.Columns.Item(1).ColumnWidth = 8 - 5/6
$ws.Columns.Item(2).ColumnWidth = 14 - 5/6
$ws.Columns.Item(3).ColumnWidth = 20 - 5/6
$ws.Columns.Item(4).ColumnWidth = 15 - 5/6
